$d = $word.ActiveDocument

# 1. Merge the two runs "Авто-" / "прокат;" into a single run "Авто-прокат;"
$d.Content.Find.Execute("Авто-прокат;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Авто-прокат;", 2) | Out-Null

# 2. Append a new paragraph after "Ключевыми объектами..." describing the
#    branches / registration process, split across several runs (as it would
#    be if typed interactively), plus the usual "_GoBack" bookmark Word drops
#    at the last edited location.
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$newParagraphXml = @"
<w:p xmlns:w="$w">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>Филиалы нашей авто-прокатной компании располагаются в н</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>ескольких городах</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t>. Клиенту необходимо пройти быструю регистрацию</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/>
    </w:rPr>
    <w:t xml:space="preserve"> предоставив ФИО, дату рождения и своё водительское удостоверение. Аренда происходит на определённое количество времени с почасовой оплатой.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@

$insertionPoint.InsertXML($newParagraphXml) | Out-Null
